$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Егор Барсуков) - mark tasks C..L, O, R as solved (1) and
# update the computed primary/secondary score columns.
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 1
$ws.Range("O2").Value = 1
$ws.Range("R2").Value = 1
$ws.Range("AD2").Value = 12
$ws.Range("AE2").Value = 56

# Row 3 (Фёдор Самохин) - clear previously solved tasks C, D, L and
# reset the computed primary/secondary score columns.
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0
